$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Header row updates
$ws.Range("B1").Value = "age.in.years"
$ws.Range("C1").Value = "Nationalität"
$ws.Range("D1").Value = "incEur"

# Row 2 - Anna -> Anna Maria
$ws.Range("A2").Value = "Anna Maria"
$ws.Range("B2").Value = 33
Set-TextValue $ws.Range("D2") "1689"

# Row 3 - Ben
$ws.Range("C3").Value = "N/A"
Set-TextValue $ws.Range("D3") "875"

# Row 4 - Clara
$ws.Range("B4").Value = 21

# Row 5 - Dimitri
Set-TextValue $ws.Range("D5") "2299"

# Row 6 - Emilia-Luise
$ws.Range("B6").Value = 29
Set-TextValue $ws.Range("D6") "2522"

# Row 7 - Fatima
$ws.Range("B7").Value = 23
Set-TextValue $ws.Range("D7") "1060"

# Row 8 - Gerda Maria -> Gerda
$ws.Range("A8").Value = "Gerda"
$ws.Range("B8").Value = 32
Set-TextValue $ws.Range("D8") "1781"

# Row 9 - Hannah
$ws.Range("B9").Value = 23
Set-TextValue $ws.Range("D9") "2463"

# Row 10 - Ismail
$ws.Range("B10").Value = 27
Set-TextValue $ws.Range("D10") "1442"

# Row 11 - Johanna
$ws.Range("B11").Value = 21
Set-TextValue $ws.Range("D11") "1404"
